$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 540, shifting existing rows 540:626 down to 541:627
$ws.Rows("540:540").Insert()

# Populate the newly inserted row 540 with its data
$ws.Cells.Item(540, 1).Value = 3
$ws.Cells.Item(540, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(540, 3).Value = "Coquimbo"
$ws.Cells.Item(540, 4).Value = 45218
$ws.Cells.Item(540, 5).Value = 5
$ws.Cells.Item(540, 6).Value = 100112012
$ws.Cells.Item(540, 7).Value = "Espinaca"
$ws.Cells.Item(540, 8).Value = "Sin especificar"
$ws.Cells.Item(540, 9).Value = "Primera"
$ws.Cells.Item(540, 10).Value = 125
$ws.Cells.Item(540, 11).Value = 4000
$ws.Cells.Item(540, 12).Value = 4000
$ws.Cells.Item(540, 13).Value = 4000
$ws.Cells.Item(540, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(540, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(540, 16).Value = 1333
$ws.Cells.Item(540, 17).Value = 3
$ws.Cells.Item(540, 18).Value = "Hortaliza"
